$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - copy formatting from the existing "sum" header (G1)
# so it reuses the same bold/centered/bordered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values: 0 for rows 2-11, 1 for row 12 (the most recent/save row)
$ws.Range("H2:H11").Value = 0
$ws.Range("H12").Value = 1
